$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2023-10-21 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-22 Sunday", 2)

# Replacement values for the 100 arithmetic-problem table cells, in
# row-major order (row 1 col 1..5, row 2 col 1..5, ...).
$newValues = @(
    "27+44=",
    "37+2=",
    "64-3=",
    "5+11=",
    "45-35=",
    "46+51=",
    "43+20=",
    "91-42=",
    "52+33=",
    "59+35=",
    "25+13=",
    "69+20=",
    "18-15=",
    "21-19=",
    "74-12=",
    "76-31=",
    "13+85=",
    "27+51=",
    "5+24=",
    "1+45=",
    "24+1=",
    "16+82=",
    "6+42=",
    "83-38=",
    "26-14=",
    "31-23=",
    "42+43=",
    "76-39=",
    "25+64=",
    "69-64=",
    "61+20=",
    "42-34=",
    "14+44=",
    "65-26=",
    "60+9=",
    "15+46=",
    "4+65=",
    "59+7=",
    "85-80=",
    "9+82=",
    "54+26=",
    "22+10=",
    "77-54=",
    "40-21=",
    "1+90=",
    "58-45=",
    "77-25=",
    "82-81=",
    "73-48=",
    "85-62=",
    "41-26=",
    "1+53=",
    "6+84=",
    "47+38=",
    "20+29=",
    "43+34=",
    "96-79=",
    "97-50=",
    "14+27=",
    "18+71=",
    "28-3=",
    "86-71=",
    "9+20=",
    "33+55=",
    "31+14=",
    "16+56=",
    "28-11=",
    "10+82=",
    "89-65=",
    "12+83=",
    "26+46=",
    "3+48=",
    "9+89=",
    "77-3=",
    "89+10=",
    "81-50=",
    "65+10=",
    "59-12=",
    "20-14=",
    "53-5=",
    "69-0=",
    "15+9=",
    "52+40=",
    "52-45=",
    "89-57=",
    "9-5=",
    "37+20=",
    "8+57=",
    "52+22=",
    "80-45=",
    "53+18=",
    "92-22=",
    "93-37=",
    "63-32=",
    "15+60=",
    "57-2=",
    "39+56=",
    "35-16=",
    "56+16=",
    "60-41="
)

$t = $d.Tables.Item(1)
$numRows = $t.Rows.Count
$numCols = 5
$i = 0
for ($r = 1; $r -le $numRows; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}

Write-Output "Updated $i cells"
